$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2423.1052
$ws.Range("I62").Value = 2226.7144
$ws.Range("J62").Value = 2973
$ws.Range("K62").Value = 2226.7144
$ws.Range("L62").Value = 2973
$ws.Range("M62").Value = -1602.7144
$ws.Range("N62").Value = -4221
$ws.Range("H65").Value = 2423.1052
$ws.Range("I65").Value = 2226.7144
$ws.Range("J65").Value = 2973
$ws.Range("K65").Value = 11133.572
$ws.Range("L65").Value = 14865
$ws.Range("M65").Value = -8013.572
$ws.Range("N65").Value = -21105
$ws.Range("H137").Value = 2660.6042
$ws.Range("I137").Value = 2599.65
$ws.Range("J137").Value = 2965.375
$ws.Range("K137").Value = 7798.950000000001
$ws.Range("L137").Value = 8896.125
$ws.Range("M137").Value = -5248.950000000001
$ws.Range("N137").Value = -13996.125
$ws.Range("H138").Value = 2290.9187
$ws.Range("I138").Value = 805.5789
$ws.Range("J138").Value = 5210.3794
$ws.Range("K138").Value = 2416.7367
$ws.Range("L138").Value = 15631.1382
$ws.Range("M138").Value = 2723.2633
$ws.Range("N138").Value = -25911.1382

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1292.6976
$ws.Range("I45").Value = 1048.8125
$ws.Range("J45").Value = 2002.1818
$ws.Range("K45").Value = 1048.8125
$ws.Range("L45").Value = 2002.1818
$ws.Range("M45").Value = -671.8125
$ws.Range("N45").Value = -2756.1818
$ws.Range("H74").Value = 900
$ws.Range("I74").Value = 711.5789
$ws.Range("J74").Value = 1795
$ws.Range("K74").Value = 711.5789
$ws.Range("L74").Value = 1795
$ws.Range("M74").Value = 162.4211
$ws.Range("N74").Value = -3543
$ws.Range("H77").Value = 900
$ws.Range("I77").Value = 711.5789
$ws.Range("J77").Value = 1795
$ws.Range("K77").Value = 3557.8945
$ws.Range("L77").Value = 8975
$ws.Range("M77").Value = 810.1055000000001
$ws.Range("N77").Value = -17711
$ws.Range("H122").Value = 2585.1516
$ws.Range("I122").Value = 1992.4
$ws.Range("J122").Value = 4437.5
$ws.Range("K122").Value = 5977.200000000001
$ws.Range("L122").Value = 13312.5
$ws.Range("M122").Value = -3527.200000000001
$ws.Range("N122").Value = -18212.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 976.7059
$ws.Range("I94").Value = 869.5
$ws.Range("K94").Value = 869.5
$ws.Range("M94").Value = -418.5
$ws.Range("H102").Value = 13389
$ws.Range("I102").Value = 13389
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 13389
$ws.Range("L102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -10144
$ws.Range("H133").Value = 20753.334
$ws.Range("J133").Value = 20753.334
$ws.Range("L133").Value = 20753.334
$ws.Range("N133").Value = -30873.334
$ws.Range("H134").Value = 2192.7144
$ws.Range("I134").Value = 1380.3889
$ws.Range("J134").Value = 7066.6665
$ws.Range("K134").Value = 4141.1667
$ws.Range("L134").Value = 21199.9995
$ws.Range("M134").Value = -1606.1667
$ws.Range("N134").Value = -26269.9995

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 9617938
$ws.Range("I58").Value = 1615.2632
$ws.Range("J58").Value = 35719384
$ws.Range("K58").Value = 1615.2632
$ws.Range("L58").Value = 35719384
$ws.Range("M58").Value = -1412.2632
$ws.Range("N58").Value = -35719790
$ws.Range("H86").Value = 6632.4165
$ws.Range("I86").Value = 5650
$ws.Range("K86").Value = 5650
$ws.Range("M86").Value = -4527
$ws.Range("H89").Value = 6632.4165
$ws.Range("I89").Value = 5650
$ws.Range("K89").Value = 28250
$ws.Range("M89").Value = -22634
$ws.Range("H122").Value = 4045.5625
$ws.Range("I122").Value = 3126.9092
$ws.Range("J122").Value = 6066.6
$ws.Range("K122").Value = 9380.7276
$ws.Range("L122").Value = 18199.8
$ws.Range("M122").Value = -6930.7276
$ws.Range("N122").Value = -23099.8
$ws.Range("H127").Value = 35000
$ws.Range("J127").Value = 35000
$ws.Range("L127").Value = 35000
$ws.Range("N127").Value = -44920
$ws.Range("H132").Value = 1687.8727
$ws.Range("I132").Value = 1324.5385
$ws.Range("J132").Value = 2573.5
$ws.Range("K132").Value = 3973.6155
$ws.Range("L132").Value = 7720.5
$ws.Range("M132").Value = -1443.6155
$ws.Range("N132").Value = -12780.5
$ws.Range("H133").Value = 28000
$ws.Range("J133").Value = 28000
$ws.Range("L133").Value = 28000
$ws.Range("N133").Value = -33060
$ws.Range("H134").Value = 1757.9269
$ws.Range("I134").Value = 1019.5
$ws.Range("J134").Value = 5344.5713
$ws.Range("K134").Value = 3058.5
$ws.Range("L134").Value = 16033.7139
$ws.Range("M134").Value = -523.5
$ws.Range("N134").Value = -21103.7139
$ws.Range("H136").Value = 9617938
$ws.Range("I136").Value = 1615.2632
$ws.Range("J136").Value = 35719384
$ws.Range("K136").Value = 4845.7896
$ws.Range("L136").Value = 107158152
$ws.Range("M136").Value = -2295.7896
$ws.Range("N136").Value = -107163252

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 770.03705
$ws.Range("I5").Value = 335.95456
$ws.Range("J5").Value = 2680
$ws.Range("K5").Value = 1007.86368
$ws.Range("L5").Value = 8040
$ws.Range("M5").Value = -895.86368
$ws.Range("N5").Value = -8264
$ws.Range("H12").Value = 203.27586
$ws.Range("I12").Value = 36.25
$ws.Range("J12").Value = 266.90475
$ws.Range("K12").Value = 108.75
$ws.Range("L12").Value = 800.71425
$ws.Range("M12").Value = 64.25
$ws.Range("N12").Value = -1146.71425
$ws.Range("H131").Value = 2139.8845
$ws.Range("I131").Value = 3147.5
$ws.Range("J131").Value = 1692.0555
$ws.Range("K131").Value = 9442.5
$ws.Range("L131").Value = 5076.166499999999
$ws.Range("M131").Value = -4402.5
$ws.Range("N131").Value = -15156.1665
$ws.Range("H135").Value = 770.03705
$ws.Range("I135").Value = 335.95456
$ws.Range("J135").Value = 2680
$ws.Range("K135").Value = 3023.59104
$ws.Range("L135").Value = 24120
$ws.Range("M135").Value = -488.5910400000002
$ws.Range("N135").Value = -29190

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H108").Value = 28442
$ws.Range("J108").Value = 28442
$ws.Range("L108").Value = 28442
$ws.Range("N108").Value = -36122
$ws.Range("H109").Value = 10256.667
$ws.Range("J109").Value = 10256.667
$ws.Range("L109").Value = 10256.667
$ws.Range("N109").Value = -12336.667

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2133.75
$ws.Range("I100").Value = 1300
$ws.Range("J100").Value = 2411.6667
$ws.Range("K100").Value = 1300
$ws.Range("L100").Value = 2411.6667
$ws.Range("M100").Value = -759
$ws.Range("N100").Value = -3493.6667
$ws.Range("H132").Value = 2628.5806
$ws.Range("I132").Value = 1613
$ws.Range("K132").Value = 4839
$ws.Range("M132").Value = -2309
$ws.Range("H136").Value = 2355.0784
$ws.Range("I136").Value = 1288.7858
$ws.Range("K136").Value = 3866.3574
$ws.Range("M136").Value = -1316.3574

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 35643.08
$ws.Range("J16").Value = 35643.08
$ws.Range("L16").Value = 35643.08
$ws.Range("N16").Value = -36227.08
$ws.Range("H81").Value = 997.1429000000001
$ws.Range("I81").Value = 997.1429000000001
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 1994.2858
$ws.Range("L81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -933.2858000000001
$ws.Range("H84").Value = 997.1429000000001
$ws.Range("I84").Value = 997.1429000000001
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 9971.429
$ws.Range("L84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -4667.429
$ws.Range("H107").Value = 564.0303
$ws.Range("I107").Value = 247.625
$ws.Range("J107").Value = 1407.7778
$ws.Range("K107").Value = 742.875
$ws.Range("L107").Value = 4223.3334
$ws.Range("M107").Value = 1177.125
$ws.Range("N107").Value = -8063.3334
$ws.Range("H132").Value = 12475.42
$ws.Range("I132").Value = 2621.2058
$ws.Range("J132").Value = 33415.625
$ws.Range("K132").Value = 7863.617400000001
$ws.Range("L132").Value = 100246.875
$ws.Range("M132").Value = -5333.617400000001
$ws.Range("N132").Value = -105306.875
$ws.Range("H136").Value = 829.9216
$ws.Range("I136").Value = 395.9
$ws.Range("J136").Value = 2408.182
$ws.Range("K136").Value = 1187.7
$ws.Range("L136").Value = 7224.545999999999
$ws.Range("M136").Value = 1362.3
$ws.Range("N136").Value = -12324.546
